# Update "paises" workbook: refresh timestamp, swap San Marino/Guatemala order,
# and update their statistics plus Panama's statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 05:22"

# 2) Swap the country names for rows 106 and 107 so that Guatemala now
#    appears before San Marino (row 106 becomes Guatemala, row 107 becomes
#    San Marino), matching the new shared-string ordering.
$ws.Range("A106").Value = "Guatemala"
$ws.Range("A107").Value = "San Marino"

# 3) Update Panama's row (row 49) with refreshed statistics.
$ws.Range("B49").Value = 6200
$ws.Range("C49").Value = 179
$ws.Range("D49").Value = 484
$ws.Range("E49").Value = 5540
$ws.Range("F49").Value = 89
$ws.Range("G49").Value = 9
$ws.Range("H49").Value = 176

# 4) Update the row that now represents Guatemala (row 106) with its new data.
$ws.Range("B106").Value = 557
$ws.Range("C106").Value = 27
$ws.Range("D106").Value = 62
$ws.Range("E106").Value = 479
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 16

# 5) Update the row that now represents San Marino (row 107) with its new data.
$ws.Range("B107").Value = 553
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 64
$ws.Range("E107").Value = 448
$ws.Range("F107").Value = 5
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 41
